# The three picture placeholders already on the slide (Maunder butterfly
# diagram, Hertzsprung-Russell diagram, Moseley atomic-weight graph) each
# get a caption textbox added above them.
#
# Shapes.AddTextbox expects Left/Top/Width/Height in points, while the
# target layout is specified in EMU (914400 EMU = 1 inch = 72 pt), so each
# value below is converted with /914400*72.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$tb1 = $s.Shapes.AddTextbox(1, 629377/914400*72, 1219200/914400*72, 4104548/914400*72, 338554/914400*72)
$tb1.Name = "TextBox 10"
$tb1.Fill.Visible = $false
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1
$tb1.TextFrame.TextRange.Text = "Maunder, butterfly diagram"
$tb1.TextFrame.TextRange.Font.Size = 16
$tb1.TextFrame.TextRange.LanguageID = "en-CA"

$tb2 = $s.Shapes.AddTextbox(1, 5304707/914400*72, 1219200/914400*72, 2153370/914400*72, 338554/914400*72)
$tb2.Name = "TextBox 11"
$tb2.Fill.Visible = $false
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1
$tb2.TextFrame.TextRange.Text = "Hertzsprung-Russell"
$tb2.TextFrame.TextRange.Font.Size = 16
$tb2.TextFrame.TextRange.LanguageID = "en-CA"

$tb3 = $s.Shapes.AddTextbox(1, 7861937/914400*72, 1219200/914400*72, 2329813/914400*72, 338554/914400*72)
$tb3.Name = "TextBox 12"
$tb3.Fill.Visible = $false
$tb3.TextFrame.WordWrap = $true
$tb3.TextFrame.AutoSize = 1
$tb3.TextFrame.TextRange.Text = "Moseley, atomic weight"
$tb3.TextFrame.TextRange.Font.Size = 16
$tb3.TextFrame.TextRange.LanguageID = "en-CA"
